$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 63.18739966666666
$ws.Range("H2").Value = 189.562199
$ws.Range("I2").Value = 0.09596345243430386
$ws.Range("J2").Value = 0.09988075390087989
$ws.Range("M2").Value = 13.89934866666667
$ws.Range("N2").Value = 41.69804600000001
$ws.Range("O2").Value = 0.04853507553134179
$ws.Range("P2").Value = 0.04999273878390351
$ws.Range("Q2").Value = 878.2636993070172
$ws.Range("R2").Value = 7904.373293763155
$ws.Range("S2").Value = 0.004657593412147263
$ws.Range("T2").Value = 0.004993312439306039
$ws.Range("G3").Value = 63.18739966666666
$ws.Range("H3").Value = 189.562199
$ws.Range("I3").Value = 0.09596345243430386
$ws.Range("J3").Value = 0.09988075390087989
$ws.Range("O3").Value = 0.245697991654417
$ws.Range("P3").Value = 0.253077086664408
$ws.Range("Q3").Value = 4446.014036248216
$ws.Range("R3").Value = 40014.12632623394
$ws.Range("S3").Value = 0.02357802753533263
$ws.Range("T3").Value = 0.02527753021107939
$ws.Range("G4").Value = 63.18739966666666
$ws.Range("H4").Value = 189.562199
$ws.Range("I4").Value = 0.09596345243430386
$ws.Range("J4").Value = 0.09988075390087989
$ws.Range("M4").Value = 82.007665
$ws.Range("N4").Value = 246.022995
$ws.Range("O4").Value = 0.2863622109480123
$ws.Range("P4").Value = 0.2949625822722868
$ws.Range("Q4").Value = 5181.851104085112
$ws.Range("R4").Value = 46636.659936766
$ws.Range("S4").Value = 0.02748030640929167
$ws.Range("T4").Value = 0.02946108508990631
$ws.Range("G5").Value = 63.18739966666666
$ws.Range("H5").Value = 189.562199
$ws.Range("I5").Value = 0.09596345243430386
$ws.Range("J5").Value = 0.09988075390087989
$ws.Range("M5").Value = 25.0501465
$ws.Range("N5").Value = 50.100293
$ws.Range("O5").Value = 0.0874724982879541
$ws.Range("P5").Value = 0.06006638442832619
$ws.Range("Q5").Value = 1582.853618604051
$ws.Range("R5").Value = 9497.121711624306
$ws.Range("S5").Value = 0.00839416292876581
$ws.Range("T5").Value = 0.005999475760801292
$ws.Range("G6").Value = 63.18739966666666
$ws.Range("H6").Value = 189.562199
$ws.Range("I6").Value = 0.09596345243430386
$ws.Range("J6").Value = 0.09988075390087989
$ws.Range("M6").Value = 95.05788666666668
$ws.Range("N6").Value = 285.17366
$ws.Range("O6").Value = 0.3319322235782747
$ws.Range("P6").Value = 0.3419012078510756
$ws.Range("Q6").Value = 6006.460676275371
$ws.Range("R6").Value = 54058.14608647834
$ws.Range("S6").Value = 0.03185336214876647
$ws.Range("T6").Value = 0.03414935039978687
$ws.Range("I7").Value = 0.3063997713314046
$ws.Range("J7").Value = 0.3189072441572365
$ws.Range("M7").Value = 13.89934866666667
$ws.Range("N7").Value = 41.69804600000001
$ws.Range("O7").Value = 0.04853507553134179
$ws.Range("P7").Value = 0.04999273878390351
$ws.Range("Q7").Value = 2804.190447614085
$ws.Range("R7").Value = 25237.71402852677
$ws.Range("S7").Value = 0.01487113604435557
$ws.Range("T7").Value = 0.01594304655344726
$ws.Range("I8").Value = 0.3063997713314046
$ws.Range("J8").Value = 0.3189072441572365
$ws.Range("O8").Value = 0.245697991654417
$ws.Range("P8").Value = 0.253077086664408
$ws.Range("S8").Value = 0.07528180845949874
$ws.Range("T8").Value = 0.08070811626748847
$ws.Range("I9").Value = 0.3063997713314046
$ws.Range("J9").Value = 0.3189072441572365
$ws.Range("M9").Value = 82.007665
$ws.Range("N9").Value = 246.022995
$ws.Range("O9").Value = 0.2863622109480123
$ws.Range("P9").Value = 0.2949625822722868
$ws.Range("Q9").Value = 16545.02785268182
$ws.Range("R9").Value = 148905.2506741364
$ws.Range("S9").Value = 0.08774131595242643
$ws.Range("T9").Value = 0.09406570424195713
$ws.Range("I10").Value = 0.3063997713314046
$ws.Range("J10").Value = 0.3189072441572365
$ws.Range("M10").Value = 25.0501465
$ws.Range("N10").Value = 50.100293
$ws.Range("O10").Value = 0.0874724982879541
$ws.Range("P10").Value = 0.06006638442832619
$ws.Range("Q10").Value = 5053.861386691842
$ws.Range("R10").Value = 30323.16832015105
$ws.Range("S10").Value = 0.02680155347321582
$ws.Range("T10").Value = 0.01915560512452665
$ws.Range("I11").Value = 0.3063997713314046
$ws.Range("J11").Value = 0.3189072441572365
$ws.Range("M11").Value = 95.05788666666668
$ws.Range("N11").Value = 285.17366
$ws.Range("O11").Value = 0.3319322235782747
$ws.Range("P11").Value = 0.3419012078510756
$ws.Range("Q11").Value = 19177.90712023165
$ws.Range("R11").Value = 172601.1640820849
$ws.Range("S11").Value = 0.101703957401908
$ws.Range("T11").Value = 0.1090347719698171
$ws.Range("G12").Value = 170.2928416666667
$ws.Range("H12").Value = 510.878525
$ws.Range("I12").Value = 0.2586257560429799
$ws.Range("J12").Value = 0.2691830570543736
$ws.Range("M12").Value = 13.89934866666667
$ws.Range("N12").Value = 41.69804600000001
$ws.Range("O12").Value = 0.04853507553134179
$ws.Range("P12").Value = 0.04999273878390351
$ws.Range("Q12").Value = 2366.959581762461
$ws.Range("R12").Value = 21302.63623586215
$ws.Range("S12").Value = 0.01255242060389641
$ws.Range("T12").Value = 0.01345719825637189
$ws.Range("G13").Value = 170.2928416666667
$ws.Range("H13").Value = 510.878525
$ws.Range("I13").Value = 0.2586257560429799
$ws.Range("J13").Value = 0.2691830570543736
$ws.Range("O13").Value = 0.245697991654417
$ws.Range("P13").Value = 0.253077086664408
$ws.Range("Q13").Value = 11982.20481166599
$ws.Range("R13").Value = 107839.8433049939
$ws.Range("S13").Value = 0.06354382884986538
$ws.Range("T13").Value = 0.06812406385873998
$ws.Range("G14").Value = 170.2928416666667
$ws.Range("H14").Value = 510.878525
$ws.Range("I14").Value = 0.2586257560429799
$ws.Range("J14").Value = 0.2691830570543736
$ws.Range("M14").Value = 82.007665
$ws.Range("N14").Value = 246.022995
$ws.Range("O14").Value = 0.2863622109480123
$ws.Range("P14").Value = 0.2949625822722868
$ws.Range("Q14").Value = 13965.31831129804
$ws.Range("R14").Value = 125687.8648016824
$ws.Range("S14").Value = 0.074060643308569
$ws.Range("T14").Value = 0.07939892961270632
$ws.Range("G15").Value = 170.2928416666667
$ws.Range("H15").Value = 510.878525
$ws.Range("I15").Value = 0.2586257560429799
$ws.Range("J15").Value = 0.2691830570543736
$ws.Range("M15").Value = 25.0501465
$ws.Range("N15").Value = 50.100293
$ws.Range("O15").Value = 0.0874724982879541
$ws.Range("P15").Value = 0.06006638442832619
$ws.Range("Q15").Value = 4265.860631651304
$ws.Range("R15").Value = 25595.16378990783
$ws.Range("S15").Value = 0.0226226410026904
$ws.Range("T15").Value = 0.01616885298662006
$ws.Range("G16").Value = 170.2928416666667
$ws.Range("H16").Value = 510.878525
$ws.Range("I16").Value = 0.2586257560429799
$ws.Range("J16").Value = 0.2691830570543736
$ws.Range("M16").Value = 95.05788666666668
$ws.Range("N16").Value = 285.17366
$ws.Range("O16").Value = 0.3319322235782747
$ws.Range("P16").Value = 0.3419012078510756
$ws.Range("Q16").Value = 16187.67764329461
$ws.Range("R16").Value = 145689.0987896515
$ws.Range("S16").Value = 0.08584622227795874
$ws.Range("T16").Value = 0.09203401233993533
$ws.Range("G17").Value = 77.473122
$ws.Range("H17").Value = 154.946244
$ws.Range("I17").Value = 0.1176593481802354
$ws.Range("J17").Value = 0.08164152846121862
$ws.Range("M17").Value = 13.89934866666667
$ws.Range("N17").Value = 41.69804600000001
$ws.Range("O17").Value = 0.04853507553134179
$ws.Range("P17").Value = 0.04999273878390351
$ws.Range("Q17").Value = 1076.825934973204
$ws.Range("R17").Value = 6460.955609839225
$ws.Range("S17").Value = 0.005710605350896164
$ws.Range("T17").Value = 0.004081483606280327
$ws.Range("G18").Value = 77.473122
$ws.Range("H18").Value = 154.946244
$ws.Range("I18").Value = 0.1176593481802354
$ws.Range("J18").Value = 0.08164152846121862
$ws.Range("O18").Value = 0.245697991654417
$ws.Range("P18").Value = 0.253077086664408
$ws.Range("Q18").Value = 5451.191054878571
$ws.Range("R18").Value = 32707.14632927142
$ws.Range("S18").Value = 0.02890866554725161
$ws.Range("T18").Value = 0.02066160017379456
$ws.Range("G19").Value = 77.473122
$ws.Range("H19").Value = 154.946244
$ws.Range("I19").Value = 0.1176593481802354
$ws.Range("J19").Value = 0.08164152846121862
$ws.Range("M19").Value = 82.007665
$ws.Range("N19").Value = 246.022995
$ws.Range("O19").Value = 0.2863622109480123
$ws.Range("P19").Value = 0.2949625822722868
$ws.Range("Q19").Value = 6353.38983548013
$ws.Range("R19").Value = 38120.33901288078
$ws.Range("S19").Value = 0.03369319108359418
$ws.Range("T19").Value = 0.02408119605557744
$ws.Range("G20").Value = 77.473122
$ws.Range("H20").Value = 154.946244
$ws.Range("I20").Value = 0.1176593481802354
$ws.Range("J20").Value = 0.08164152846121862
$ws.Range("M20").Value = 25.0501465
$ws.Range("N20").Value = 50.100293
$ws.Range("O20").Value = 0.0874724982879541
$ws.Range("P20").Value = 0.06006638442832619
$ws.Range("Q20").Value = 1940.713055912373
$ws.Range("R20").Value = 7762.852223649493
$ws.Range("S20").Value = 0.01029195713225743
$ws.Range("T20").Value = 0.004903911433867691
$ws.Range("G21").Value = 77.473122
$ws.Range("H21").Value = 154.946244
$ws.Range("I21").Value = 0.1176593481802354
$ws.Range("J21").Value = 0.08164152846121862
$ws.Range("M21").Value = 95.05788666666668
$ws.Range("N21").Value = 285.17366
$ws.Range("O21").Value = 0.3319322235782747
$ws.Range("P21").Value = 0.3419012078510756
$ws.Range("Q21").Value = 7364.431250788841
$ws.Range("R21").Value = 44186.58750473305
$ws.Range("S21").Value = 0.03905492906623595
$ws.Range("T21").Value = 0.02791333719169862
$ws.Range("G22").Value = 145.7496183333334
$ws.Range("H22").Value = 437.248855
$ws.Range("I22").Value = 0.2213516720110761
$ws.Range("J22").Value = 0.2303874164262914
$ws.Range("M22").Value = 13.89934866666667
$ws.Range("N22").Value = 41.69804600000001
$ws.Range("O22").Value = 0.04853507553134179
$ws.Range("P22").Value = 0.04999273878390351
$ws.Range("Q22").Value = 2025.824763248593
$ws.Range("R22").Value = 18232.42286923733
$ws.Range("S22").Value = 0.01074332012004637
$ws.Range("T22").Value = 0.01151769792849799
$ws.Range("G23").Value = 145.7496183333334
$ws.Range("H23").Value = 437.248855
$ws.Range("I23").Value = 0.2213516720110761
$ws.Range("J23").Value = 0.2303874164262914
$ws.Range("O23").Value = 0.245697991654417
$ws.Range("P23").Value = 0.253077086664408
$ws.Range("Q23").Value = 10255.28590045245
$ws.Range("R23").Value = 92297.57310407203
$ws.Range("S23").Value = 0.05438566126246863
$ws.Range("T23").Value = 0.0583057761533056
$ws.Range("G24").Value = 145.7496183333334
$ws.Range("H24").Value = 437.248855
$ws.Range("I24").Value = 0.2213516720110761
$ws.Range("J24").Value = 0.2303874164262914
$ws.Range("M24").Value = 82.007665
$ws.Range("N24").Value = 246.022995
$ws.Range("O24").Value = 0.2863622109480123
$ws.Range("P24").Value = 0.2949625822722868
$ws.Range("Q24").Value = 11952.58587415786
$ws.Range("R24").Value = 107573.2728674207
$ws.Range("S24").Value = 0.063386754194131
$ws.Range("T24").Value = 0.06795566727213957
$ws.Range("G25").Value = 145.7496183333334
$ws.Range("H25").Value = 437.248855
$ws.Range("I25").Value = 0.2213516720110761
$ws.Range("J25").Value = 0.2303874164262914
$ws.Range("M25").Value = 25.0501465
$ws.Range("N25").Value = 50.100293
$ws.Range("O25").Value = 0.0874724982879541
$ws.Range("P25").Value = 0.06006638442832619
$ws.Range("Q25").Value = 3651.049291569087
$ws.Range("R25").Value = 21906.29574941452
$ws.Range("S25").Value = 0.01936218375102463
$ws.Range("T25").Value = 0.01383853912251049
$ws.Range("G26").Value = 145.7496183333334
$ws.Range("H26").Value = 437.248855
$ws.Range("I26").Value = 0.2213516720110761
$ws.Range("J26").Value = 0.2303874164262914
$ws.Range("M26").Value = 95.05788666666668
$ws.Range("N26").Value = 285.17366
$ws.Range("O26").Value = 0.3319322235782747
$ws.Range("P26").Value = 0.3419012078510756
$ws.Range("Q26").Value = 13854.65070123993
$ws.Range("R26").Value = 124691.8563111593
$ws.Range("S26").Value = 0.07347375268340543
$ws.Range("T26").Value = 0.07876973594983777
